$d = $word.ActiveDocument

function Replace-Text($old, $new) {
    $d.Content.Find.Execute($old, $true, $false, $false, $false, $false, `
                             $true, 1, $false, $new, 2)
}

Replace-Text "2023-12-29 Friday" "2023-12-30 Saturday"

Replace-Text "88×56=4928" "63×54=3402"
Replace-Text "86×91=7826" "80×16=1280"
Replace-Text "97×46=4462" "34×63=2142"
Replace-Text "25×33=825" "77×19=1463"
Replace-Text "23×78=1794" "26×86=2236"

Replace-Text "30×37=1110" "66×17=1122"
Replace-Text "66×23=1518" "70×88=6160"
Replace-Text "25×69=1725" "32×16=512"
Replace-Text "21×12=252" "27×92=2484"
Replace-Text "33×75=2475" "46×34=1564"

Replace-Text "85×59=5015" "69×96=6624"
Replace-Text "79×39=3081" "90×62=5580"
Replace-Text "22×84=1848" "74×90=6660"
Replace-Text "94×29=2726" "81×42=3402"
Replace-Text "57×62=3534" "63×89=5607"

Replace-Text "90×80=7200" "96×63=6048"
Replace-Text "94×97=9118" "28×16=448"
Replace-Text "14×57=798" "99×54=5346"
Replace-Text "28×23=644" "64×61=3904"
Replace-Text "18×46=828" "69×73=5037"

Replace-Text "29×50=1450" "78×99=7722"
Replace-Text "45×26=1170" "62×38=2356"
Replace-Text "66×75=4950" "37×93=3441"
Replace-Text "47×98=4606" "53×85=4505"
Replace-Text "57×26=1482" "93×56=5208"
